# Swap the data of row 4 <-> row 5, and row 6 <-> row 7
# (Columns A,B,D,E,F,G,H,Q,R plus the presence of an (empty) AF cell.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    $cols = @("A","B","D","E","F","G","H","Q","R")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }

    # Swap presence of the (empty) "Bestämningsmetod" cell in column AF:
    # whichever of the two rows originally had the (blank) AF cell should
    # have it after the swap, and vice versa. Touching NumberFormat (without
    # changing it) is enough to make an otherwise-empty cell persist in the
    # saved sheet, which is how we recreate a "present but blank" cell.
    $had1 = $ws.Range("AF$r1").Value2 -ne $null
    $had2 = $ws.Range("AF$r2").Value2 -ne $null
    $ws.Range("AF$r1").ClearContents()
    $ws.Range("AF$r2").ClearContents()
    if ($had2) {
        $ws.Range("AF$r1").NumberFormat = $ws.Range("AF$r1").NumberFormat
    }
    if ($had1) {
        $ws.Range("AF$r2").NumberFormat = $ws.Range("AF$r2").NumberFormat
    }
}

Swap-RowData 4 5
Swap-RowData 6 7
